$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(32)
Write-Host "text=[$($p.Range.Text)]"
$xml = $p.Range.WordOpenXML
Write-Host "xml=[$xml]"
